$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scattered cell corrections in rows 2-25 ---

# Row 2
$ws.Range("D2").Value = ""

# Row 4
$ws.Range("E4").Value = ""

# Row 5
$ws.Range("D5").Value = -14.4

# Row 6
$ws.Range("C6").Value = 15.1
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7

# Row 8
$ws.Range("C8").Value = ""

# Row 9
$ws.Range("D9").Value = ""

# Row 10
$ws.Range("D10").Value = ""

# Row 11
$ws.Range("E11").Value = -7.9

# Row 12
$ws.Range("C12").Value = 12.5
$ws.Range("E12").Value = ""

# Row 14
$ws.Range("C14").Value = ""
$ws.Range("E14").Value = -5.4

# Row 17
$ws.Range("C17").Value = 11.2
$ws.Range("E17").Value = ""

# Row 18
$ws.Range("C18").Value = 11.5

# Row 19
$ws.Range("C19").Value = ""
$ws.Range("E19").Value = -6.5

# Row 20
$ws.Range("C20").Value = ""

# Row 21
$ws.Range("E21").Value = -8.699999999999999

# Row 22
$ws.Range("E22").Value = -6.1

# Row 23
$ws.Range("C23").Value = 12.2

# Row 24
$ws.Range("D24").Value = -13.9

# Row 25
$ws.Range("E25").Value = ""

# --- Rows 26-33 replaced with the next set of records (IDs/values shift up) ---

$newRows = @(
  @("SC 5",   -20.2, 10.8,  -13.8, "",    17.38),
  @("SC 101", -20.4, "",    -14.6, "",    17),
  @("SC 105", "",    11.1,  "",    "",    17.44),
  @("SC 119", "",    11.2,  -13,   -6.8,  18.06),
  @("SC 120", -19.7, "",    -13.6, -5.7,  16.89),
  @("SC 132", -18.8, 15.3,  -13.7, -8.1,  17.18),
  @("SC 193", "",    10.5,  -14.7, -6.4,  17.39),
  @("SC 232", -19.5, 10.4,  -14.1, -10.7, 17.53)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowNum = 26 + $i
    $rowData = $newRows[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

# --- Rows 34 and 35 no longer exist; delete them entirely ---
$ws.Range("A34:F35").Delete() | Out-Null
